# Update average_county_temperature (column I) values with refreshed NOAA data,
# and recompute the dependent worst/best ASHP COP values (columns N/O) for the
# rows that had them populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# average_county_temperature (column I) updates
$ws.Range("I2").Value  = 1.925925925925943
$ws.Range("I9").Value  = 1.925925925925943
$ws.Range("I10").Value = 1.925925925925943
$ws.Range("I11").Value = 1.925925925925943
$ws.Range("I12").Value = -1.226851851851833
$ws.Range("I14").Value = 1.925925925925943
$ws.Range("I18").Value = 13.17361111111111
$ws.Range("I19").Value = 13.17361111111111
$ws.Range("I20").Value = 13.17361111111111

# worst_ashp_cop (column N) / best_ashp_cop (column O) recomputed for rows
# whose heat pump COP depends on the updated temperature
$ws.Range("N10").Value = 1.204711500590784
$ws.Range("O10").Value = 1.258324667221298

$ws.Range("N19").Value = 1.281341554412284
$ws.Range("O19").Value = 1.343162802314449
